$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("SFIA Level"), shifting the
# existing SFIA Level / Keycode / Description columns one to the right.
$ws.Columns.Item(2).Insert()

# New header for the inserted column.
$ws.Range("B1").Value = "Skill Description"

# Map each SkillCode (column A) to its full skill name for the new
# "Skill Description" column (column B).
$skillNames = @{
    "Autonomy" = "Autonomy";
    "Influence" = "Influence";
    "Complexity" = "Complexity";
    "Knowledge" = "Knowledge";
    "SCTY" = "Information security";
    "CNSL" = "Consultancy";
    "BPRE" = "Business process improvement";
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

for ($r = 2; $r -le $lastRow; $r++) {
    $code = $ws.Cells.Item($r, 1).Value2
    if ($code -and $skillNames.ContainsKey($code)) {
        $ws.Cells.Item($r, 2).Value = $skillNames[$code]
    }
}
